# Weekly update: insert a new price record as row 39, pushing the existing
# rows 39-109 down to 40-110 (see commit message: "Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 39; this shifts rows 39-109
# down to 40-110, carrying all of their data/formatting along with them.
$ws.Rows("39:39").Insert()

# Populate the newly inserted row 39 with the new weekly data point.
# (Same market/category/variety/quality/unit/origin as the former row 39,
#  now at row 40, but a newer date and new price figures.)
$ws.Range("A39").Value = 11
$ws.Range("B39").Value = "Vega Monumental Concepción"
$ws.Range("C39").Value = "Bíobío"
$ws.Range("D39").Value = 44868
$ws.Range("D39").NumberFormat = $ws.Range("D40").NumberFormat
$ws.Range("E39").Value = 8
$ws.Range("F39").Value = 100112024
$ws.Range("G39").Value = "Choclo"
$ws.Range("H39").Value = "Dulce o Americano"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 100
$ws.Range("K39").Value = 25000
$ws.Range("L39").Value = 26000
$ws.Range("M39").Value = 25500
$ws.Range("N39").Value = '$/malla 70 unidades'
$ws.Range("O39").Value = "Región de Arica y Parinacota"
$ws.Range("P39").Value = 364
$ws.Range("Q39").Value = 70
$ws.Range("R39").Value = "Hortaliza"
